# Rescale detector radius measurements to cm.
#
# Columns N (minR) and O (maxR) on Sheet2 were entered in millimetres;
# convert every data row (4-33) to centimetres by dividing by 10, matching
# the units already used elsewhere (column C etc. use [cm]).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

$firstRow = 4
$lastRow  = 33
$colMinR  = 14   # N
$colMaxR  = 15   # O

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $minRCell = $ws.Cells.Item($r, $colMinR)
    $maxRCell = $ws.Cells.Item($r, $colMaxR)

    $minRCell.Value = $minRCell.Value2 / 10
    $maxRCell.Value = $maxRCell.Value2 / 10
}

# Drop the stale O4:O33 selection that was left over from inspecting the
# old millimetre values.
$ws.Range("A1").Select()
